# 010 Week 6 Pairs data update
# Fill in WK 6 (column I) scores on Sheet1 for the pairs competition.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$wk6 = @{
    12 = 32
    14 = 33
    15 = 30
    16 = 28
    17 = 34
    18 = 36
    20 = 38
    21 = 31
    22 = 25
    23 = 27
    25 = 26
    26 = 30
    27 = 21
    28 = 29
    29 = 36
    31 = 39
}

foreach ($row in $wk6.Keys) {
    $ws.Range("I$row").Value = $wk6[$row]
}
